$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Row 2
$ws.Range("C2").Value = 11.22
$ws.Range("E2").Value = 11.12

# Row 3
$ws.Range("B3").Value = 8.779999999999999
$ws.Range("D3").Value = 10.42
$ws.Range("E3").Value = 10.55
$ws.Range("F3").Value = 9.31

# Row 4
$ws.Range("C4").Value = 9.58
$ws.Range("E4").Value = 10.46
$ws.Range("F4").Value = 10.21

# Row 5
$ws.Range("B5").Value = 8.880000000000001
$ws.Range("C5").Value = 9.449999999999999
$ws.Range("D5").Value = 9.539999999999999
$ws.Range("F5").Value = 10.29
$ws.Range("G5").Value = 9.34
$ws.Range("H5").Value = 8.27

# Row 6
$ws.Range("C6").Value = 10.69
$ws.Range("D6").Value = 9.789999999999999
$ws.Range("E6").Value = 9.710000000000001
$ws.Range("G6").Value = 10.42
$ws.Range("H6").Value = 10.9

# Row 7
$ws.Range("E7").Value = 10.66
$ws.Range("F7").Value = 9.58

# Row 8
$ws.Range("E8").Value = 11.73
$ws.Range("F8").Value = 9.1
